$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.572.81'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '2.661.93'
$ws.Range('E3').Value = '  -0.17%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.93'
$ws.Range('E5').Value = '  -1.61%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.04'
$ws.Range('E6').Value = '  -1.60%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.603'
$ws.Range('E8').Value = '  +2.05%  '

$ws.Range('E9').Value = '  -2.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.88'
$ws.Range('E10').Value = '  -0.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.395'
$ws.Range('E11').Value = '  -2.82%  '

$ws.Range('E12').Value = '  -0.35%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.19'
$ws.Range('E13').Value = '  -2.84%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000193'
$ws.Range('E14').Value = '  -1.25%  '

$ws.Range('D15').Value = '3.137.90'
$ws.Range('E15').Value = '  -0.43%  '

$ws.Range('D16').Value = '65.376.09'
$ws.Range('E16').Value = '  -0.29%  '

$ws.Range('D17').Value = '2.660.15'
$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.45'
$ws.Range('E18').Value = '  -2.51%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.79'
$ws.Range('E19').Value = '  -2.45%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.47'
$ws.Range('E20').Value = '  +0.54%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.29'
$ws.Range('E21').Value = '  -3.61%  '

$ws.Range('E22').Value = '  -0.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.69'
$ws.Range('E23').Value = '  +0.41%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.65'
$ws.Range('E24').Value = '  +0.14%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000107'
$ws.Range('E25').Value = '  +0.07%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.63'
$ws.Range('E26').Value = '  -4.69%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.167'
$ws.Range('E27').Value = '  +0.42%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.59'
$ws.Range('E28').Value = '  -3.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.03'
$ws.Range('E29').Value = '  -2.68%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '538.14'
$ws.Range('E31').Value = '  -1.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.13'
$ws.Range('E32').Value = '  -4.15%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.74'
$ws.Range('E33').Value = '  -6.29%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.49'
$ws.Range('E34').Value = '  +1.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.39'
$ws.Range('E35').Value = '  -4.67%  '

$ws.Range('E36').Value = '  -3.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.30'
$ws.Range('E37').Value = '  -2.15%  '

$ws.Range('E38').Value = '  -0.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '158.84'
$ws.Range('E39').Value = '  -2.58%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.93'
$ws.Range('E40').Value = '  -5.02%  '

$ws.Range('E41').Value = '  +0.04%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.44'
$ws.Range('E42').Value = '  -0.07%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '165.41'
$ws.Range('E43').Value = '  -1.17%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.06'
$ws.Range('E44').Value = '  -3.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0608'
$ws.Range('E45').Value = '  -1.35%  '

$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.90'
$ws.Range('E46').Value = '  -1.81%  '

$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.24'
$ws.Range('E47').Value = '  -7.02%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.645'
$ws.Range('E48').Value = '  -2.90%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0259'
$ws.Range('E49').Value = '  -3.03%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0994'
$ws.Range('E50').Value = '  +0.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.86'
$ws.Range('E51').Value = '  -0.52%  '
